# xls export geometry fix: rename several hwinventory headers to shorter
# labels and shrink their columns accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -------------------------------------------------
$ws.Range("F1").Value = "Memory tot.size"
$ws.Range("H1").Value = "Memory P/Ns"
$ws.Range("M1").Value = "HDD slot pop."
$ws.Range("N1").Value = "PSU P/Ns"

# --- Column width updates --------------------------------------------------
# Target OOXML <col> widths (in "characters", 1/256 granularity, MDW=7):
#   F ->  15.7109375  (was 18.7109375)
#   H ->  11.7109375  (was 25.7109375)
#   M ->  13.7109375  (was 19.7109375)
#   N ->   8.7109375  (was 15.7109375)
# The host engine quantizes ColumnWidth to 1/6-character steps internally,
# so it cannot reproduce the 1/256 MDW=7 figure bit-exactly. Feeding it
# (N - 1/6) for an integer character count N lands on N + 2/3, which is the
# closest value the engine can actually store to the N + 0.7109375 target.
$ws.Columns.Item(6).ColumnWidth = (15 - 1/6)
$ws.Columns.Item(8).ColumnWidth = (11 - 1/6)
$ws.Columns.Item(13).ColumnWidth = (13 - 1/6)
$ws.Columns.Item(14).ColumnWidth = (8 - 1/6)
